# Refresh the cryptos price/volume table (columns D and E) on Sheet1.
# Some "Price" values parse as plain numbers (e.g. "199.41", "8.60"), so we
# force the cell to Text format first to keep them as literal strings
# (matching the source data's original formatting, e.g. trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.476.00"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "2.958.04"
$ws.Range("E3").Value = "  +2.17%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.41"
$ws.Range("E5").Value = "  +1.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "597.08"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("E9").Value = "  +4.23%  "

$ws.Range("D10").Value = "2.957.07"
$ws.Range("E10").Value = "  +2.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.445"
$ws.Range("E11").Value = "  +11.02%  "

$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").Value = "3.501.93"
$ws.Range("E14").Value = "  +2.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.63"
$ws.Range("E15").Value = "  +4.40%  "

$ws.Range("D16").Value = "76.397.28"
$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "2.950.17"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.69"
$ws.Range("E19").Value = "  +8.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.82"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.02"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("E23").Value = "  +4.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.68"
$ws.Range("E24").Value = "  +1.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.35"
$ws.Range("E27").Value = "  +2.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.73"
$ws.Range("E28").Value = "  -1.05%  "

$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.60"
$ws.Range("E31").Value = "  +10.15%  "

$ws.Range("E32").Value = "  -1.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "498.38"
$ws.Range("E33").Value = "  -2.14%  "

$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.90"
$ws.Range("E36").Value = "  +1.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.36"
$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.393"
$ws.Range("E38").Value = "  +13.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("E39").Value = "  +18.12%  "

$ws.Range("E40").Value = "  +1.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  -2.51%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "181.53"
$ws.Range("E43").Value = "  -0.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.95"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.66"
$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").Value = "  -2.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.75"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.592"
$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("E49").Value = "  +3.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  -2.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.671"
$ws.Range("E51").Value = "  +0.33%  "
